# Insert a new weekly price record before the current row 25, shifting the
# remaining historical rows down by one (dimension grows from R53 to R54).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(25).Insert()

$ws.Cells.Item(25, 1).Value = 8
$ws.Cells.Item(25, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(25, 3).Value = "Coquimbo"
$ws.Cells.Item(25, 4).Value = 44482
$ws.Cells.Item(25, 5).Value = 4
$ws.Cells.Item(25, 6).Value = 100112028
$ws.Cells.Item(25, 7).Value = "Sandia"
$ws.Cells.Item(25, 8).Value = "Sin especificar"
$ws.Cells.Item(25, 9).Value = "Primera"
$ws.Cells.Item(25, 10).Value = 800
$ws.Cells.Item(25, 11).Value = 800
$ws.Cells.Item(25, 12).Value = 900
$ws.Cells.Item(25, 13).Value = 850
$ws.Cells.Item(25, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(25, 15).Value = "Perú"
$ws.Cells.Item(25, 16).Value = 850
$ws.Cells.Item(25, 17).Value = 1
$ws.Cells.Item(25, 18).Value = "Hortaliza"
